$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Fecha) updates
$ws.Range("D2").Value = 44302
$ws.Range("D3").Value = 44322
$ws.Range("D4").Value = 44309
$ws.Range("D5").Value = 44316
$ws.Range("D6").Value = 44327
$ws.Range("D7").Value = 44306
$ws.Range("D8").Value = 44323
$ws.Range("D9").Value = 44313
$ws.Range("D10").Value = 44330

# Column M (Volumen) updates
$ws.Range("M2").Value = 80
$ws.Range("M3").Value = 60
$ws.Range("M5").Value = 120
$ws.Range("M6").Value = 60
$ws.Range("M7").Value = 80
$ws.Range("M8").Value = 80
$ws.Range("M9").Value = 120
$ws.Range("M10").Value = 60

# Row 4 Q/S/T updates
$ws.Range("Q4").Value = "$/caja 14 kilos granel"
$ws.Range("S4").Value = 821
$ws.Range("T4").Value = 14

# Row 6 Q/S/T updates
$ws.Range("Q6").Value = "$/caja 10 kilos empedrada"
$ws.Range("S6").Value = 11500
$ws.Range("T6").Value = 1
